# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / LevePrice / LeveProfit figures across sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1668.6666
$ws.Range("I28").Value = 1878
$ws.Range("J28").Value = 1250
$ws.Range("K28").Value = 1878
$ws.Range("L28").Value = 1250
$ws.Range("M28").Value = -1393
$ws.Range("N28").Value = -2220

$ws.Range("H41").Value = 899.8823
$ws.Range("I41").Value = 171.55556
$ws.Range("J41").Value = 1719.25
$ws.Range("K41").Value = 171.55556
$ws.Range("L41").Value = 1719.25
$ws.Range("M41").Value = 268.44444
$ws.Range("N41").Value = -2599.25

$ws.Range("H55").Value = 155.88889
$ws.Range("I55").Value = 160.5
$ws.Range("K55").Value = 160.5
$ws.Range("M55").Value = 53.5

$ws.Range("H129").Value = 2045.1428
$ws.Range("J129").Value = 4989.3335
$ws.Range("L129").Value = 14968.0005
$ws.Range("N129").Value = -24968.0005

$ws.Range("H138").Value = 24454.562
$ws.Range("I138").Value = 2983.7368
$ws.Range("J138").Value = 38521.656
$ws.Range("K138").Value = 8951.2104
$ws.Range("L138").Value = 115564.968
$ws.Range("M138").Value = -3811.2104
$ws.Range("N138").Value = -125844.968

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21337.873
$ws.Range("J32").Value = 1708.1
$ws.Range("L32").Value = 1708.1
$ws.Range("N32").Value = -2282.1

$ws.Range("H61").Value = 10297.143
$ws.Range("I61").Value = 6917.273
$ws.Range("K61").Value = 6917.273
$ws.Range("M61").Value = -6705.273

$ws.Range("H112").Value = 74500
$ws.Range("J112").Value = 74500
$ws.Range("L112").Value = 74500
$ws.Range("N112").Value = -77454

$ws.Range("H132").Value = 2025.36
$ws.Range("I132").Value = 1757.6842
$ws.Range("J132").Value = 2873
$ws.Range("K132").Value = 5273.0526
$ws.Range("L132").Value = 8619
$ws.Range("M132").Value = -2743.0526
$ws.Range("N132").Value = -13679

$ws.Range("H136").Value = 10297.143
$ws.Range("I136").Value = 6917.273
$ws.Range("K136").Value = 20751.819
$ws.Range("M136").Value = -18201.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 632
$ws.Range("I80").Value = 996.75
$ws.Range("K80").Value = 996.75
$ws.Range("M80").Value = 1.25

$ws.Range("H82").Value = 13676

$ws.Range("H83").Value = 632
$ws.Range("I83").Value = 996.75
$ws.Range("K83").Value = 4983.75
$ws.Range("M83").Value = 8.25

$ws.Range("H85").Value = 13676

$ws.Range("H99").Value = 2331.3333
$ws.Range("I99").Value = 2136.8667
$ws.Range("J99").Value = 3303.6667
$ws.Range("K99").Value = 2136.8667
$ws.Range("L99").Value = 3303.6667
$ws.Range("M99").Value = -638.8667
$ws.Range("N99").Value = -6299.6667

$ws.Range("H105").Value = 2582.1304
$ws.Range("I105").Value = 2231.2942
$ws.Range("J105").Value = 3576.1667
$ws.Range("K105").Value = 2231.2942
$ws.Range("L105").Value = 3576.1667
$ws.Range("M105").Value = -484.2941999999998
$ws.Range("N105").Value = -7070.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1589832.2
$ws.Range("I31").Value = 3032343.2
$ws.Range("J31").Value = 3070.2
$ws.Range("K31").Value = 3032343.2
$ws.Range("L31").Value = 3070.2
$ws.Range("M31").Value = -3032048.2
$ws.Range("N31").Value = -3660.2

$ws.Range("H34").Value = 1589832.2
$ws.Range("I34").Value = 3032343.2
$ws.Range("J34").Value = 3070.2
$ws.Range("K34").Value = 3032343.2
$ws.Range("L34").Value = 3070.2
$ws.Range("M34").Value = -3032141.2
$ws.Range("N34").Value = -3474.2

$ws.Range("H134").Value = 2146.6943
$ws.Range("I134").Value = 1955.9688
$ws.Range("J134").Value = 3672.5
$ws.Range("K134").Value = 5867.9064
$ws.Range("L134").Value = 11017.5
$ws.Range("M134").Value = -3332.9064
$ws.Range("N134").Value = -16087.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 117.6
$ws.Range("I14").Value = 117.6
$ws.Range("K14").Value = 352.8
$ws.Range("M14").Value = -179.8

$ws.Range("H36").Value = 254.5
$ws.Range("I36").Value = 254.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 763.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -594.5
$ws.Range("N36").ClearContents()

$ws.Range("H132").Value = 1761
$ws.Range("I132").Value = 1437.25
$ws.Range("J132").Value = 2020
$ws.Range("K132").Value = 12935.25
$ws.Range("L132").Value = 18180
$ws.Range("M132").Value = -10405.25
$ws.Range("N132").Value = -23240

$ws.Range("H133").Value = 2399.8
$ws.Range("I133").Value = 1749.75
$ws.Range("K133").Value = 5249.25
$ws.Range("M133").Value = -189.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 12000000
$ws.Range("I10").Value = 12000000
$ws.Range("K10").Value = 12000000
$ws.Range("M10").Value = -11999831

$ws.Range("H11").Value = 31501
$ws.Range("J11").Value = 31501
$ws.Range("L11").Value = 31501
$ws.Range("N11").Value = -31779

$ws.Range("H80").Value = 3545.4
$ws.Range("I80").Value = 2734.3333
$ws.Range("K80").Value = 2734.3333
$ws.Range("M80").Value = -1736.3333

$ws.Range("H83").Value = 3545.4
$ws.Range("I83").Value = 2734.3333
$ws.Range("K83").Value = 13671.6665
$ws.Range("M83").Value = -8679.666499999999

$ws.Range("H132").Value = 2144.3333
$ws.Range("I132").Value = 2001.6
$ws.Range("K132").Value = 6004.799999999999
$ws.Range("M132").Value = -3474.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3288.4167
$ws.Range("I61").Value = 3147.1
$ws.Range("K61").Value = 3147.1
$ws.Range("M61").Value = -2945.1

$ws.Range("H100").Value = 3333.3333
$ws.Range("J100").Value = 4000
$ws.Range("L100").Value = 4000
$ws.Range("N100").Value = -5082

$ws.Range("H113").Value = 3288.4167
$ws.Range("I113").Value = 3147.1
$ws.Range("K113").Value = 3147.1
$ws.Range("M113").Value = -977.0999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H70").Value = 95399.60000000001
$ws.Range("J70").Value = 95399.60000000001
$ws.Range("L70").Value = 95399.60000000001
$ws.Range("N70").Value = -96029.60000000001

$ws.Range("H73").Value = 95399.60000000001
$ws.Range("J73").Value = 95399.60000000001
$ws.Range("L73").Value = 95399.60000000001
$ws.Range("N73").Value = -97583.60000000001

$ws.Range("H123").Value = 138999.75
$ws.Range("J123").Value = 138999.75
$ws.Range("L123").Value = 138999.75
$ws.Range("N123").Value = -148799.75

$ws.Range("H132").Value = 22695.4
$ws.Range("I132").Value = 35217.383
$ws.Range("J132").Value = 3912.4285
$ws.Range("K132").Value = 105652.149
$ws.Range("L132").Value = 11737.2855
$ws.Range("M132").Value = -103122.149
$ws.Range("N132").Value = -16797.2855

